$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Início")
$ws2 = $wb.Worksheets.Item("Cad_Empresa")

# --- Cad_Empresa sheet: fill in the address number / city used by the
# concatenation formula on the Início sheet ---
$ws2.Range("B14").Value = "rua oscar pereira"
$ws2.Range("B22").Value = "catu"

# --- Início sheet: new G-column demo values for TEXT()/date formulas ---
$ws1.Range("G7").Value = 43520
$ws1.Range("G7").NumberFormat = "mm-dd-yy"

$ws1.Range("G8").Value = 0.28299999999999997
$ws1.Range("G9").Value = 234
$ws1.Range("G10").Value = 1245

$ws1.Range("G12").Formula = '=TEXT(G10,"R$#.##0,00")'
$ws1.Range("G14").Formula = '=TEXT(G9,"0000")'
$ws1.Range("G15").Formula = '=TEXT(G8,"0,00%")'
$ws1.Range("G16").Formula = '=TEXT(G7,"dd/mm/aaaa")'
$ws1.Range("G17").Formula = '=TEXT(G7,"dd/mm/aa")'
$ws1.Range("G18").Formula = '=TEXT(G7,"dddd")'
$ws1.Range("G19").Formula = '=TEXT(G7,"mmmm")'

foreach ($addr in @("G12","G14","G15","G16","G17","G18","G19")) {
    $c = $ws1.Range($addr)
    $c.Font.Name = "Tahoma"
    $c.Font.Size = 10
    $c.HorizontalAlignment = -4152
}

# Concatenation formula (nested PROPER + CONCATENATE) building the address line
$ws1.Range("B17").Formula = '=CONCATENATE(PROPER( Cad_Empresa!B14),", Número ",Cad_Empresa!B18," - ", PROPER(Cad_Empresa!B22))'

# New notes rows about the concatenation lesson
$ws1.Range("B24").Value = "Concatenar células(conteúdos) : insere o sinal = juntamente com a fórmula concatenar/ clica na célula /clica no sinal ;"
$ws1.Range("B25").Value = "seleciona a segunda célula/ quando tiver todas as células que deseja clica no enter"
$ws1.Range("B26").Value = "Editar fórmula : fn+f2"
$ws1.Range("B27").Value = "Espaço entre as concatenações: só inserir apas duplas entre os parâmetros isolados pelo sinal ;"
$ws1.Range("B28").Value = "aninhamento de fórmulas, exemplo: =CONCATENAR(PRI.MAIÚSCULA(Cad_Empresa!B14);"

# --- selection / active-cell bookkeeping (match the author's final cursor position) ---
$ws2.Activate()
$ws2.Range("B6").Select()

$ws1.Activate()
$ws1.Range("H14").Select()
